$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Stamp formatting for the 5 new rows (192-196) by copying the
#        formats (borders/font/number-format/alignment) from the last
#        existing data row (191) straight down. Format-only paste keeps the
#        existing style table untouched (reuses style ids 1 and 2).
$ws.Range("A191:V191").Copy()
$ws.Range("A192:V196").PasteSpecial(-4122)

# Column D ("temporada") is the literal text "2023" on every existing row.
# Assigning the string "2023" via .Value would be reinterpreted as a number
# by Excel's type inference, so instead copy the already-text D191 cell's
# VALUE (not format) down - a plain values-paste keeps the text type.
for ($r = 192; $r -le 196; $r++) {
    $ws.Range("D191").Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4163)
}

$excel.CutCopyMode = 0

# --- 2) Fill in the new match rows scraped on 08-11-2023.
$data = @(
    @(191,"colombia","primera-a",45238.0625,"Petrolera",2,"Pereira",1,2.18,"01/11/2023 10:42",1.98,"08/11/2023 01:22",3.1,"01/11/2023 10:42",3.24,"08/11/2023 01:20",3.57,"01/11/2023 10:42",4.51,"08/11/2023 01:22","https://www.betexplorer.com/football/colombia/primera-a/petrolera-dep-pereira/QePEbNKH/"),
    @(192,"colombia","primera-a",45238.0625,"Chico",1,"Dep. Cali",1,2.65,"01/11/2023 10:43",3.75,"08/11/2023 01:23",2.97,"01/11/2023 10:43",3.1,"08/11/2023 01:23",2.9,"01/11/2023 10:43",2.25,"08/11/2023 01:23","https://www.betexplorer.com/football/colombia/primera-a/chico-dep-cali/pGp3UN4o/"),
    @(193,"colombia","primera-a",45238.0625,"Envigado",1,"Dep. Pasto",1,2.7,"01/11/2023 10:42",2.83,"08/11/2023 01:29",2.92,"01/11/2023 10:42",3.26,"08/11/2023 01:29",3.03,"01/11/2023 10:42",2.69,"08/11/2023 01:29","https://www.betexplorer.com/football/colombia/primera-a/envigado-dep-pasto/ptUg3o6n/"),
    @(194,"colombia","primera-a",45238.0625,"Junior",2,"Huila",0,1.51,"01/11/2023 10:42",1.3,"08/11/2023 01:28",4.13,"01/11/2023 10:42",5.33,"08/11/2023 01:28",7.03,"01/11/2023 10:42",12.68,"08/11/2023 01:28","https://www.betexplorer.com/football/colombia/primera-a/junior-huila/8KEJcszO/"),
    @(195,"colombia","primera-a",45238.0625,"Santa Fe",0,"Once Caldas",1,1.82,"01/11/2023 10:43",2.29,"08/11/2023 01:27",3.43,"01/11/2023 10:43",3.33,"08/11/2023 01:27",4.45,"01/11/2023 10:43",3.36,"08/11/2023 01:27","https://www.betexplorer.com/football/colombia/primera-a/santa-fe-once-caldas/n5QAa35B/")
)

$r = 192
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]    # Indice
    $ws.Cells.Item($r, 2).Value = $row[1]    # pais
    $ws.Cells.Item($r, 3).Value = $row[2]    # torneio
    # column 4 (temporada) already populated above via values-paste
    $ws.Cells.Item($r, 5).Value = $row[3]    # data_partida
    $ws.Cells.Item($r, 6).Value = $row[4]    # home
    $ws.Cells.Item($r, 7).Value = $row[5]    # home_ft_gols
    $ws.Cells.Item($r, 8).Value = $row[6]    # away
    $ws.Cells.Item($r, 9).Value = $row[7]    # away_ft_gols
    $ws.Cells.Item($r, 10).Value = $row[8]   # home_opening_odds
    $ws.Cells.Item($r, 11).Value = $row[9]   # home_opening_data_hora
    $ws.Cells.Item($r, 12).Value = $row[10]  # home_closing_odds
    $ws.Cells.Item($r, 13).Value = $row[11]  # home_closing_data_hora
    $ws.Cells.Item($r, 14).Value = $row[12]  # draw_opening_odds
    $ws.Cells.Item($r, 15).Value = $row[13]  # draw_opening_data_hora
    $ws.Cells.Item($r, 16).Value = $row[14]  # draw_closing_odds
    $ws.Cells.Item($r, 17).Value = $row[15]  # draw_closing_data_hora
    $ws.Cells.Item($r, 18).Value = $row[16]  # away_opening_odds
    $ws.Cells.Item($r, 19).Value = $row[17]  # away_opening_data_hora
    $ws.Cells.Item($r, 20).Value = $row[18]  # away_closing_odds
    $ws.Cells.Item($r, 21).Value = $row[19]  # away_closing_data_hora
    $ws.Cells.Item($r, 22).Value = $row[20]  # url_partida
    $r = $r + 1
}

Write-Host "Added rows 192-196"
